$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = '033/2024'
$ws.Range("B32").Value = '003/2023'
$ws.Range("C32").Value = 'Não informado'
$ws.Range("D32").Value = '07 de agosto de 2024'
$ws.Range("E32").Value = 'RABELLO CONSTRUÇÕES LTDA, CNPJ 26.695.883/0001-95'
$ws.Range("F32").Value = 'Tomada de Preços'
$ws.Range("G32").Value = 'Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.'
$ws.Range("H32").Value = 'PREFEITURA MUNICIPAL DE NILO PEÇANHA – BA.'
$ws.Range("I32").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos).'

$ws.Range("A33").Value = '033/2024'
$ws.Range("B33").Value = '003/2023'
$ws.Range("C33").Value = 'Não informado'
$ws.Range("D33").Value = '07 de agosto de 2024'
$ws.Range("E33").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F33").Value = 'Tomada de Preços'
$ws.Range("G33").Value = '"seleção de proposta mais vantajosa para Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.", consoante às condições estabelecidas no edital da Tomada de Preços Nº 003/2023.'
$ws.Range("H33").Value = 'MUNICÍPIO DE NILO PEÇANHA, Estado da Bahia'
$ws.Range("I33").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos).'

$ws.Range("A34").Value = '033/2024'
$ws.Range("B34").Value = '003/2023'
$ws.Range("C34").Value = 'Não informado'
$ws.Range("D34").Value = '07 de agosto de 2024'
$ws.Range("E34").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F34").Value = 'Tomada de Preços'
$ws.Range("G34").Value = '"seleção de proposta mais vantajosa para Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.", consoante às condições estabelecidas no edital da Tomada de Preços Nº 003/2023.'
$ws.Range("H34").Value = 'MUNICÍPIO DE NILO PEÇANHA, CNPJ/MF N° 13.758.313/0001-55'
$ws.Range("I34").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos)'

$ws.Range("A35").Value = '033/2024'
$ws.Range("B35").Value = '003/2023'
$ws.Range("C35").Value = 'Não informado'
$ws.Range("D35").Value = '07 de agosto de 2024'
$ws.Range("E35").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F35").Value = 'Tomada de Preços'
$ws.Range("G35").Value = 'Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.'
$ws.Range("H35").Value = 'MUNICÍPIO DE NILO PEÇANHA, CNPJ/MF N° 13.758.313/0001-55'
$ws.Range("I35").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos)'

$ws.Range("A36").Value = '033/2024'
$ws.Range("B36").Value = '003/2023'
$ws.Range("C36").Value = 'Não informado'
$ws.Range("D36").Value = '07 de agosto de 2024'
$ws.Range("E36").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F36").Value = 'Tomada de Preços'
$ws.Range("G36").Value = 'Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.'
$ws.Range("H36").Value = 'MUNICÍPIO DE NILO PEÇANHA'
$ws.Range("I36").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos)'

$ws.Range("A37").Value = '033/2024'
$ws.Range("B37").Value = '003/2023'
$ws.Range("C37").Value = 'Não informado'
$ws.Range("D37").Value = '07 de agosto de 2024'
$ws.Range("E37").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F37").Value = 'Tomada de Preços'
$ws.Range("G37").Value = 'Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.'
$ws.Range("H37").Value = 'MUNICÍPIO DE NILO PEÇANHA, CNPJ/MF N° 13.758.313/0001-55'
$ws.Range("I37").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos)'

$ws.Range("A38").Value = '033/2024'
$ws.Range("B38").Value = '003/2023'
$ws.Range("C38").Value = 'Não informado'
$ws.Range("D38").Value = '07 de agosto de 2024'
$ws.Range("E38").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F38").Value = 'Tomada de Preços'
$ws.Range("G38").Value = '"seleção de proposta mais vantajosa para Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.", consoante às condições estabelecidas no edital da Tomada de Preços Nº 003/2023.'
$ws.Range("H38").Value = 'MUNICÍPIO DE NILO PEÇANHA, CNPJ/MF N° 13.758.313/0001-55'
$ws.Range("I38").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos).'

$ws.Range("A39").Value = '033/2024'
$ws.Range("B39").Value = '003/2023'
$ws.Range("C39").Value = 'Não informado'
$ws.Range("D39").Value = '07 de agosto de 2024'
$ws.Range("E39").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F39").Value = 'Tomada de Preços'
$ws.Range("G39").Value = 'Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.'
$ws.Range("H39").Value = 'MUNICÍPIO DE NILO PEÇANHA, CNPJ/MF N° 13.758.313/0001-55'
$ws.Range("I39").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos)'

$ws.Range("A40").Value = '033/2024'
$ws.Range("B40").Value = '003/2023'
$ws.Range("C40").Value = 'Não informado'
$ws.Range("D40").Value = '07 de agosto de 2024'
$ws.Range("E40").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F40").Value = 'Tomada de Preços'
$ws.Range("G40").Value = 'Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.'
$ws.Range("H40").Value = 'MUNICÍPIO DE NILO PEÇANHA, CNPJ/MF N° 13.758.313/0001-55'
$ws.Range("I40").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos).'

$ws.Range("A41").Value = '033/2024'
$ws.Range("B41").Value = '003/2023'
$ws.Range("C41").Value = 'Não informado'
$ws.Range("D41").Value = '07 de agosto de 2024'
$ws.Range("E41").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F41").Value = 'Tomada de Preços'
$ws.Range("G41").Value = 'Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.'
$ws.Range("H41").Value = 'MUNICÍPIO DE NILO PEÇANHA, Estado da Bahia'
$ws.Range("I41").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos)'

$ws.Range("A42").Value = '033/2024'
$ws.Range("B42").Value = '003/2023'
$ws.Range("C42").Value = 'Não informado'
$ws.Range("D42").Value = '07 de agosto de 2024'
$ws.Range("E42").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F42").Value = 'Tomada de Preços'
$ws.Range("G42").Value = 'Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.'
$ws.Range("H42").Value = 'MUNICÍPIO DE NILO PEÇANHA, Estado da Bahia'
$ws.Range("I42").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos).'

$ws.Range("A43").Value = '033/2024'
$ws.Range("B43").Value = '003/2023'
$ws.Range("C43").Value = 'Não informado'
$ws.Range("D43").Value = '07 de agosto de 2024'
$ws.Range("E43").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F43").Value = 'Tomada de Preços'
$ws.Range("G43").Value = 'Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.'
$ws.Range("H43").Value = 'MUNICÍPIO DE NILO PEÇANHA, Estado da Bahia'
$ws.Range("I43").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos)'

$ws.Range("A44").Value = '033/2024'
$ws.Range("B44").Value = '003/2023'
$ws.Range("C44").Value = 'Não informado'
$ws.Range("D44").Value = '07 de agosto de 2024'
$ws.Range("E44").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F44").Value = 'Tomada de Preços'
$ws.Range("G44").Value = 'Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.'
$ws.Range("H44").Value = 'MUNICÍPIO DE NILO PEÇANHA'
$ws.Range("I44").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos).'

$ws.Range("A45").Value = '033/2024'
$ws.Range("B45").Value = '003/2023'
$ws.Range("C45").Value = 'Não informado'
$ws.Range("D45").Value = '07 de agosto de 2024'
$ws.Range("E45").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F45").Value = 'Tomada de Preços'
$ws.Range("G45").Value = '"seleção de proposta mais vantajosa para Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.", consoante às condições estabelecidas no edital da Tomada de Preços Nº 003/2023.'
$ws.Range("H45").Value = 'MUNICÍPIO DE NILO PEÇANHA, Estado da Bahia'
$ws.Range("I45").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos).'

$ws.Range("A46").Value = '033/2024'
$ws.Range("B46").Value = '003/2023'
$ws.Range("C46").Value = 'Não informado'
$ws.Range("D46").Value = '07 de agosto de 2024'
$ws.Range("E46").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F46").Value = 'Tomada de Preços'
$ws.Range("G46").Value = 'Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.'
$ws.Range("H46").Value = 'MUNICÍPIO DE NILO PEÇANHA, Estado da Bahia'
$ws.Range("I46").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos)'

$ws.Range("A47").Value = '033/2024'
$ws.Range("B47").Value = '003/2023'
$ws.Range("C47").Value = 'Não informado'
$ws.Range("D47").Value = '07 de agosto de 2024'
$ws.Range("E47").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F47").Value = 'Tomada de Preços'
$ws.Range("G47").Value = 'Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.'
$ws.Range("H47").Value = 'MUNICÍPIO DE NILO PEÇANHA, Estado da Bahia'
$ws.Range("I47").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos)'

$ws.Range("A48").Value = '033/2024'
$ws.Range("B48").Value = '003/2023'
$ws.Range("C48").Value = 'Não informado'
$ws.Range("D48").Value = '07 de agosto de 2024'
$ws.Range("E48").Value = 'RABELLO CONSTRUÇÕES LTDA., CNPJ 26.695.883/0001-95'
$ws.Range("F48").Value = 'Tomada de Preços'
$ws.Range("G48").Value = 'Execução da obra da segunda etapa da requalificação da orla fluvial do Rio das Almas, na Sede do Município de Nilo Peçanha BA.'
$ws.Range("H48").Value = 'MUNICÍPIO DE NILO PEÇANHA, Estado da Bahia, pessoa jurídica de direito público interno, inscrita no CNPJ/MF N° 13.758.313/0001-55, situada na Rua Raimundo Brito, s/nº, centro, СЕР. 45.440-000, na cidade de Nilo Peçanha-Ba., fone (73) 3257-2434'
$ws.Range("I48").Value = 'R$ 637.919,93 (Seiscentos e trinta e sete mil, novecentos e dezanove reais e noventa e três centavos)'
